$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "294.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.38%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.85%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.949"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.23%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07333"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.96%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.824"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-13.91%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.667"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.78%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.754"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.81%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9081"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.84%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1657"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.07%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07606"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.02%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08193"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-7.07%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02982"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.28%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1000"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.05%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001506"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.56%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005680"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.07%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.464"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.28%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-7.60%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.07%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.30%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.372"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.17%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04480"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.59%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001226"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.74%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004048"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-10.00%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001252"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.40%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01666"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.99%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04408"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-7.29%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007410"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.96%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1324"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.04%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002057"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.65%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01116"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.99%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005966"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.55%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.17%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.126"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "158.18%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.17%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.17%"
